$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C9) from serial 45208 (2023-10-09)
# to serial 45212 (2023-10-13), preserving existing cell formatting.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45212
}
